# =====================================================================
# Mother-In-Law House Expenses - Balance reconciliation & fee cleanup
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value to a cell while forcing Excel to keep it
# as text (prevents "18.59%" / "79.8%" style strings from being
# silently re-interpreted as numeric percentages).
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# =====================================================================
# Sheet: Home Summary
# =====================================================================
$ws1 = $wb.Worksheets.Item("Home Summary")

$ws1.Range("B5").Value = "KES 798,000"
$ws1.Range("B6").Value = "KES 202,000"
Set-TextValue $ws1.Range("B7") "79.8%"
$ws1.Range("B8").Value = "KES 8,909"

$ws1.Range("B11").Value = "KES 30,000"
$ws1.Range("B13").Value = "KES 30,000"
$ws1.Range("B14").Value = "KES 828,000"
Set-TextValue $ws1.Range("B15") "82.80%"
$ws1.Range("B16").Value = "KES 172,000"

# Category Breakdown table (rows 20-28) - values reshuffled and updated
$ws1.Range("A20").Value = "Building Materials"
$ws1.Range("B20").Value = 183580
$ws1.Range("C20").Value = 2300
$ws1.Range("D20").Value = 185880
Set-TextValue $ws1.Range("E20") "18.59%"

$ws1.Range("A21").Value = "Metal & Steel"
$ws1.Range("B21").Value = 176310
$ws1.Range("C21").Value = 2090.5
$ws1.Range("D21").Value = 178400.5
Set-TextValue $ws1.Range("E21") "17.84%"

$ws1.Range("A22").Value = "Labor Costs"
$ws1.Range("B22").Value = 172700
$ws1.Range("C22").Value = 1972.5
$ws1.Range("D22").Value = 174672.5
Set-TextValue $ws1.Range("E22") "17.47%"

# Row 23 (Hardware Items) is unchanged.

$ws1.Range("A24").Value = "Miscellaneous"
$ws1.Range("B24").Value = 78545.5
$ws1.Range("C24").Value = 470
$ws1.Range("D24").Value = 79015.5
Set-TextValue $ws1.Range("E24") "7.90%"

$ws1.Range("A25").Value = "Workers Accommodation"
$ws1.Range("B25").Value = 17760
$ws1.Range("C25").Value = 190
$ws1.Range("D25").Value = 17950
Set-TextValue $ws1.Range("E25") "1.80%"

$ws1.Range("A26").Value = "Transport & Logistics"
$ws1.Range("B26").Value = 16550
$ws1.Range("C26").Value = 185
$ws1.Range("D26").Value = 16735
Set-TextValue $ws1.Range("E26") "1.67%"

$ws1.Range("A27").Value = "Utilities"
$ws1.Range("B27").Value = 8330
$ws1.Range("C27").Value = 85
$ws1.Range("D27").Value = 8415
Set-TextValue $ws1.Range("E27") "0.84%"

$ws1.Range("A28").Value = "Utilities & Services"
$ws1.Range("B28").Value = 5100
$ws1.Range("C28").Value = 75
$ws1.Range("D28").Value = 5175
Set-TextValue $ws1.Range("E28") "0.52%"

# =====================================================================
# Sheet: Daily Expenses
# =====================================================================
$ws2 = $wb.Worksheets.Item("Daily Expenses")

# Row 247: was a "Labor Payment" lump sum -> becomes a personal-use
# ("Self") miscellaneous withdrawal, reconciled to the real bank balance.
$ws2.Range("B247").Value = "Miscellaneous"
$ws2.Range("C247").Value = "Not used for the project"
$ws2.Range("D247").Value = "Money that has been taken out for personal use."
$ws2.Range("E247").Value = 77045.5
$ws2.Range("G247").Value = 77500.5
$ws2.Range("H247").Value = "Self"

# Insert a new row 248 (shifts old row 248 "Transaction Fees" down to 249)
$ws2.Rows.Item(248).Insert()
$ws2.Range("A247:I247").Copy()
$ws2.Range("A248:I248").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("A248").Value = "24/09/2025"
$ws2.Range("B248").Value = "Labor Costs"
$ws2.Range("C248").Value = "Specialized Labor"
$ws2.Range("D248").Value = "Welder - final balance payment (18000 of 40000 total)"
$ws2.Range("E248").Value = 18000
$ws2.Range("F248").Value = 171.5
$ws2.Range("G248").Value = 18171.5
$ws2.Range("H248").Value = "Welder"
$ws2.Range("I248").Value = "PAID"

# Old row 248 ("Transaction Fees", now shifted to row 249) is replaced
# entirely by the Excavation Worker final balance payment.
$ws2.Range("A249").Value = "25/09/2025"
$ws2.Range("B249").Value = "Labor Costs"
$ws2.Range("C249").Value = "Specialized Labor"
$ws2.Range("D249").Value = "Excavation Worker - final balance payment (6400 of 16000 total)"
$ws2.Range("E249").Value = 6400
$ws2.Range("F249").Value = 75
$ws2.Range("G249").Value = 6475
$ws2.Range("H249").Value = "Excavation Worker"
$ws2.Range("I249").Value = "PAID"

# =====================================================================
# Sheet: M-Pesa Fees
# =====================================================================
$ws4 = $wb.Worksheets.Item("M-Pesa Fees")

$ws4.Range("C7").Value = 11
$ws4.Range("D7").Value = 825

$ws4.Range("C9").Value = 5
$ws4.Range("D9").Value = 857.5

$ws4.Range("C15").Value = 43
$ws4.Range("D15").Value = 215

$ws4.Range("B20").Value = "KES 8,909"

# =====================================================================
# Sheet: Outstanding Balances
# =====================================================================
$ws5 = $wb.Worksheets.Item("Outstanding Balances")

# Welder and Excavation Worker balances are now fully paid (see Daily
# Expenses rows 248/249) so their outstanding-balance rows are removed;
# remaining rows (Electrician, Plumber) shift up, as does the total.
$ws5.Range("A4:D5").Delete(-4162)   # xlShiftUp

$ws5.Range("C7").Value = "KES 30,000"

Write-Host "All edits applied."
